# Applies the edit described in the commit:
#   "Fix typo and removed 2035_TM152_DBP_NoProject_04"
#
# Concretely, on worksheet "all_runs":
#   1. Delete the entire row that holds the "2035_TM152_DBP_NoProject_04" run,
#      shifting all following rows up by one.
#   2. Fix the typo "RTP2022" -> "RTP2021" (found in column A of the last
#      data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("all_runs")

# 1. Locate and delete the row for the removed run "2035_TM152_DBP_NoProject_04".
$runCell = $ws.Columns.Item(3).Find("2035_TM152_DBP_NoProject_04")
if ($runCell -ne $null) {
    $ws.Rows.Item($runCell.Row).Delete()
}

# 2. Fix the "RTP2022" typo -> "RTP2021".
$typoCell = $ws.Columns.Item(1).Find("RTP2022")
if ($typoCell -ne $null) {
    $typoCell.Value = "RTP2021"
}

# Keep the selection consistent with the edited sheet (last data row, col A).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$ws.Cells.Item($lastRow, 1).Select()
